$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2185
$ws.Range("I111").Value = 2333.3333
$ws.Range("J111").Value = 2096
$ws.Range("K111").Value = 6999.999899999999
$ws.Range("L111").Value = 6288
$ws.Range("M111").Value = -3932.999899999999
$ws.Range("N111").Value = -12422

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4501.273
$ws.Range("J112").Value = 1406.4
$ws.Range("L112").Value = 4219.200000000001
$ws.Range("N112").Value = -6435.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 692.82355
$ws.Range("I132").Value = 692.82355
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2078.47065
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 451.5293500000002
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1320.9254
$ws.Range("I45").Value = 1209.2222
$ws.Range("J45").Value = 1784.9231
$ws.Range("K45").Value = 1209.2222
$ws.Range("L45").Value = 1784.9231
$ws.Range("M45").Value = -832.2221999999999
$ws.Range("N45").Value = -2538.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4923.5137
$ws.Range("I74").Value = 2376.4211
$ws.Range("J74").Value = 7612.1113
$ws.Range("K74").Value = 2376.4211
$ws.Range("L74").Value = 7612.1113
$ws.Range("M74").Value = -1502.4211
$ws.Range("N74").Value = -9360.1113

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4923.5137
$ws.Range("I77").Value = 2376.4211
$ws.Range("J77").Value = 7612.1113
$ws.Range("K77").Value = 11882.1055
$ws.Range("L77").Value = 38060.5565
$ws.Range("M77").Value = -7514.1055
$ws.Range("N77").Value = -46796.5565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 56498.332
$ws.Range("J137").Value = 56498.332
$ws.Range("L137").Value = 56498.332
$ws.Range("N137").Value = -66698.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 40000
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 10780
$ws.Range("J58").Value = 10780
$ws.Range("L58").Value = 10780
$ws.Range("N58").Value = -11368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 69400
$ws.Range("J124").Value = 69400
$ws.Range("L124").Value = 69400
$ws.Range("N124").Value = -79220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 40000
$ws.Range("N136").Value = -50200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1897995
$ws.Range("J58").Value = 7569.4546
$ws.Range("L58").Value = 7569.4546
$ws.Range("N58").Value = -7975.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3432.4443
$ws.Range("I132").Value = 3417.0908
$ws.Range("K132").Value = 10251.2724
$ws.Range("M132").Value = -7721.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1897995
$ws.Range("J136").Value = 7569.4546
$ws.Range("L136").Value = 22708.3638
$ws.Range("N136").Value = -27808.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4211.9565
$ws.Range("I3").Value = 3079
$ws.Range("J3").Value = 5083.4614
$ws.Range("K3").Value = 9237
$ws.Range("L3").Value = 15250.3842
$ws.Range("M3").Value = -9125
$ws.Range("N3").Value = -15474.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 38759.543
$ws.Range("I14").Value = 38759.543
$ws.Range("K14").Value = 116278.629
$ws.Range("M14").Value = -116105.629

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4144.875
$ws.Range("I63").Value = 992
$ws.Range("J63").Value = 4595.2856
$ws.Range("K63").Value = 2976
$ws.Range("L63").Value = 13785.8568
$ws.Range("M63").Value = -2227
$ws.Range("N63").Value = -15283.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 4144.875
$ws.Range("I66").Value = 992
$ws.Range("J66").Value = 4595.2856
$ws.Range("K66").Value = 8928
$ws.Range("L66").Value = 41357.5704
$ws.Range("M66").Value = -5184
$ws.Range("N66").Value = -48845.5704

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2668.7231
$ws.Range("J68").Value = 3590.9285
$ws.Range("L68").Value = 10772.7855
$ws.Range("N68").Value = -12394.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2668.7231
$ws.Range("J71").Value = 3590.9285
$ws.Range("L71").Value = 32318.3565
$ws.Range("N71").Value = -40430.3565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3125
$ws.Range("J103").Value = 4000
$ws.Range("L103").Value = 12000
$ws.Range("N103").Value = -13758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 616.2857
$ws.Range("I113").Value = 591.6818
$ws.Range("J113").Value = 706.5
$ws.Range("K113").Value = 1775.0454
$ws.Range("L113").Value = 2119.5
$ws.Range("M113").Value = 394.9546
$ws.Range("N113").Value = -6459.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2658.7727
$ws.Range("I126").Value = 1780
$ws.Range("K126").Value = 5340
$ws.Range("M126").Value = -2870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2799.3333
$ws.Range("I7").Value = 2118.8
$ws.Range("K7").Value = 2118.8
$ws.Range("M7").Value = -2006.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -20814

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 11000
$ws.Range("J48").Value = 18000
$ws.Range("L48").Value = 18000
$ws.Range("N48").Value = -19322

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 66000
$ws.Range("J124").Value = 66000
$ws.Range("L124").Value = 66000
$ws.Range("N124").Value = -75820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2799.3333
$ws.Range("I126").Value = 2118.8
$ws.Range("K126").Value = 6356.400000000001
$ws.Range("M126").Value = -3886.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4715.6177
$ws.Range("J136").Value = 5118.3335
$ws.Range("L136").Value = 15355.0005
$ws.Range("N136").Value = -20455.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 78800
$ws.Range("J109").Value = 78800
$ws.Range("L109").Value = 78800
$ws.Range("N109").Value = -81574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 835.63635
$ws.Range("I113").Value = 456
$ws.Range("J113").Value = 947.2941
$ws.Range("K113").Value = 1368
$ws.Range("L113").Value = 2841.8823
$ws.Range("M113").Value = 802
$ws.Range("N113").Value = -7181.882299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 58533.332
$ws.Range("J120").Value = 58533.332
$ws.Range("L120").Value = 58533.332
$ws.Range("N120").Value = -68209.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2251.5
$ws.Range("I122").Value = 1802.4
$ws.Range("K122").Value = 5407.200000000001
$ws.Range("M122").Value = -2957.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1556.9286
$ws.Range("I126").Value = 1602
$ws.Range("J126").Value = 1496.8334
$ws.Range("K126").Value = 4806
$ws.Range("L126").Value = 4490.5002
$ws.Range("M126").Value = -2336
$ws.Range("N126").Value = -9430.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5876.85
$ws.Range("I136").Value = 2453.9285
$ws.Range("J136").Value = 7719.9614
$ws.Range("K136").Value = 7361.7855
$ws.Range("L136").Value = 23159.8842
$ws.Range("M136").Value = -4811.7855
$ws.Range("N136").Value = -28259.8842
